$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Actual" results (and derived "Correct" flag) for games
#     that have since been played (rows 58-61) ---
$ws.Range("F58").Value = "Calgary Hitmen"
$ws.Range("G58").Value = 1

$ws.Range("F59").Value = "Lethbridge Hurricanes"
$ws.Range("G59").Value = 0

$ws.Range("F60").Value = "Wenatchee Wild"
$ws.Range("G60").Value = 0

$ws.Range("F61").Value = "Spokane Chiefs"
$ws.Range("G61").Value = 1

# --- Append the next day's slate of games (rows 62-65) ---
$ws.Range("A62").Value = 1021608
$ws.Range("B62").Value = "Tue, Jan 7, 2025"
$ws.Range("C62").Value = "Lethbridge Hurricanes"
$ws.Range("D62").Value = "Moose Jaw Warriors"
$ws.Range("E62").Value = "Moose Jaw Warriors"

$ws.Range("A63").Value = 1021610
$ws.Range("B63").Value = "Tue, Jan 7, 2025"
$ws.Range("C63").Value = "Red Deer Rebels"
$ws.Range("D63").Value = "Saskatoon Blades"
$ws.Range("E63").Value = "Saskatoon Blades"

$ws.Range("A64").Value = 1021609
$ws.Range("B64").Value = "Tue, Jan 7, 2025"
$ws.Range("C64").Value = "Prince George Cougars"
$ws.Range("D64").Value = "Vancouver Giants"
$ws.Range("E64").Value = "Vancouver Giants"

$ws.Range("A65").Value = 1021611
$ws.Range("B65").Value = "Tue, Jan 7, 2025"
$ws.Range("C65").Value = "Seattle Thunderbirds"
$ws.Range("D65").Value = "Victoria Royals"
$ws.Range("E65").Value = "Victoria Royals"

# --- Update the selection / active cell to reflect the newly added rows ---
$ws.Range("A62:E65").Select()
